# Remove duplicate for "Carbon dioxide, non-fossil" / "air::non-urban air or
# from high stacks" (data row 8, i.e. worksheet row 8). Deleting the entire
# row shifts every following row up by one, which is exactly what the target
# workbook shows (dimension shrinks from C224 to C223 and every row's data
# below the duplicate moves up by one position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("8").Delete()

# Restore a plain, unscrolled view with the selection on A14 (matches the
# saved sheetView/selection in the edited workbook).
$excel.Goto($ws.Range("A1"), $true)
$ws.Range("A14").Select()
